$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows (header unchanged, row2 folder_id value changes)
$ws.Range("B2").Value = "10-tourisme"

# Add new data rows (entering folder_id column before doc_id column to
# reproduce the original authoring order of new shared-string values)
$ws.Range("B3").Value = "bevnat"
$ws.Range("A3").Value = "bevnat_info"

$ws.Range("A4").Value = "bevnat_variable"
$ws.Range("B4").Value = "bevnat"

$ws.Range("B5").Value = "statpop"
$ws.Range("A5").Value = "statpop_info"

# Expand the table to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B5"))

# Move selection to A6, matching the post-edit saved state
$ws.Range("A6").Select()
